$d = $word.ActiveDocument

# 1. Remove the leading "Chad, We are going to create an annotate bibliography..."
#    paragraph entirely (the instructor's prompt that preceded the bibliography
#    entry / annotation text).
$d.Paragraphs(1).Range.Delete()

# 2. Clean up the stale lastRenderedPageBreak split in the "In conclusion" ...
#    paragraph by touching text inside the second half of the split run; this
#    causes Word to re-merge the two adjacent, identically-formatted runs
#    (dropping the now-stale rendered-page-break marker) while leaving the
#    surrounding grammar-check markers untouched.
$d.Content.Find.Execute("addresses previous", $true, $false, $false, $false, `
    $false, $true, 1, $false, "addresses previous", 2) | Out-Null
